# allow any type of SourceInput in AssemblySource
# The "input" header column (a SourceInput list) used to sit near the front
# of several *Source sheets. Move it so it sits immediately before the
# trailing "id" column on each of the AssemblySource-family sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "AssemblySource",
    "PCRSource",
    "LigationSource",
    "HomologousRecombinationSource",
    "GibsonAssemblySource",
    "InFusionSource",
    "OverlapExtensionPCRLigationSource",
    "InVivoAssemblySource",
    "RestrictionAndLigationSource",
    "GatewaySource",
    "CreLoxRecombinationSource",
    "CRISPRSource"
)

$maxCols = 10

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Collect current header values (row 1), stopping at first blank.
    $headers = @()
    for ($c = 1; $c -le $maxCols; $c++) {
        $v = $ws.Cells.Item(1, $c).Value2
        if ($v -eq $null -or $v -eq "") {
            break
        }
        $headers += $v
    }

    # Remove "input" from wherever it is, then re-insert it right before "id".
    $newHeaders = @()
    foreach ($h in $headers) {
        if ($h -ne "input") {
            $newHeaders += $h
        }
    }

    $idPos = $newHeaders.Count
    for ($i = 0; $i -lt $newHeaders.Count; $i++) {
        if ($newHeaders[$i] -eq "id") {
            $idPos = $i
        }
    }

    $finalHeaders = @()
    for ($i = 0; $i -lt $idPos; $i++) {
        $finalHeaders += $newHeaders[$i]
    }
    $finalHeaders += "input"
    for ($i = $idPos; $i -lt $newHeaders.Count; $i++) {
        $finalHeaders += $newHeaders[$i]
    }

    for ($c = 1; $c -le $finalHeaders.Count; $c++) {
        $ws.Cells.Item(1, $c).Value = $finalHeaders[$c - 1]
    }
}
